$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 224.8
$ws.Range("I12").Value = 230.5
$ws.Range("J12").Value = 202
$ws.Range("K12").Value = 230.5
$ws.Range("L12").Value = 202
$ws.Range("M12").Value = -60.5
$ws.Range("N12").Value = -542

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H41").Value = 1867
$ws.Range("I41").Value = 414
$ws.Range("K41").Value = 414
$ws.Range("M41").Value = 26

$ws.Range("H53").Value = 83334264
$ws.Range("J53").Value = 72
$ws.Range("L53").Value = 72
$ws.Range("N53").Value = -1346

$ws.Range("H125").Value = 925

$ws.Range("H126").Value = 77748.75
$ws.Range("J126").Value = 77748.75
$ws.Range("L126").Value = 77748.75
$ws.Range("N126").Value = -87628.75

$ws.Range("H127").Value = 1116.8125
$ws.Range("J127").Value = 1109.2
$ws.Range("L127").Value = 3327.6
$ws.Range("N127").Value = -13247.6

$ws.Range("H137").Value = 4212.8213
$ws.Range("J137").Value = 5064.5625
$ws.Range("L137").Value = 15193.6875
$ws.Range("N137").Value = -20293.6875

$ws.Range("H138").Value = 526822.25
$ws.Range("J138").Value = 813794.3
$ws.Range("L138").Value = 2441382.9
$ws.Range("N138").Value = -2451662.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2158.3914
$ws.Range("I2").Value = 2124.6667
$ws.Range("J2").Value = 2279.8
$ws.Range("K2").Value = 2124.6667
$ws.Range("L2").Value = 2279.8
$ws.Range("M2").Value = -2011.6667
$ws.Range("N2").Value = -2505.8

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H32").Value = 16484.195
$ws.Range("I32").Value = 10367.821
$ws.Range("K32").Value = 10367.821
$ws.Range("M32").Value = -10080.821

$ws.Range("H45").Value = 31303.625
$ws.Range("I45").Value = 56369.125
$ws.Range("K45").Value = 56369.125
$ws.Range("M45").Value = -55992.125

$ws.Range("H61").Value = 5764.2915
$ws.Range("I61").Value = 3155.4
$ws.Range("J61").Value = 10112.444
$ws.Range("K61").Value = 3155.4
$ws.Range("L61").Value = 10112.444
$ws.Range("M61").Value = -2943.4
$ws.Range("N61").Value = -10536.444

$ws.Range("H116").Value = 2158.3914
$ws.Range("I116").Value = 2124.6667
$ws.Range("J116").Value = 2279.8
$ws.Range("K116").Value = 2124.6667
$ws.Range("L116").Value = 2279.8
$ws.Range("M116").Value = 169.3332999999998
$ws.Range("N116").Value = -6867.8

$ws.Range("H136").Value = 5764.2915
$ws.Range("I136").Value = 3155.4
$ws.Range("J136").Value = 10112.444
$ws.Range("K136").Value = 9466.200000000001
$ws.Range("L136").Value = 30337.332
$ws.Range("M136").Value = -6916.200000000001
$ws.Range("N136").Value = -35437.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2158.3914
$ws.Range("I3").Value = 2124.6667
$ws.Range("J3").Value = 2279.8
$ws.Range("K3").Value = 2124.6667
$ws.Range("L3").Value = 2279.8
$ws.Range("M3").Value = -2010.6667
$ws.Range("N3").Value = -2507.8

$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H134").Value = 3568.25
$ws.Range("I134").Value = 3125.8235
$ws.Range("K134").Value = 9377.470499999999
$ws.Range("M134").Value = -6842.470499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2537.8333
$ws.Range("I16").Value = 1045.5
$ws.Range("K16").Value = 1045.5
$ws.Range("M16").Value = -758.5

$ws.Range("H62").Value = 14297500
$ws.Range("J62").Value = 14500
$ws.Range("L62").Value = 14500
$ws.Range("N62").Value = -15748

$ws.Range("H65").Value = 14297500
$ws.Range("J65").Value = 14500
$ws.Range("L65").Value = 72500
$ws.Range("N65").Value = -78740

$ws.Range("H113").Value = 2537.8333
$ws.Range("I113").Value = 1045.5
$ws.Range("K113").Value = 1045.5
$ws.Range("M113").Value = 1124.5

$ws.Range("H122").Value = 5469.231
$ws.Range("I122").Value = 4922.3335
$ws.Range("J122").Value = 6699.75
$ws.Range("K122").Value = 14767.0005
$ws.Range("L122").Value = 20099.25
$ws.Range("M122").Value = -12317.0005
$ws.Range("N122").Value = -24999.25

$ws.Range("H141").Value = 901157.5600000001
$ws.Range("J141").Value = 901157.5600000001
$ws.Range("L141").Value = 901157.5600000001
$ws.Range("N141").Value = -911517.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 410.375
$ws.Range("J98").Value = 323.33334
$ws.Range("L98").Value = 970.0000200000001
$ws.Range("N98").Value = -3966.00002

$ws.Range("H113").Value = 7222.0557
$ws.Range("J113").Value = 7222.0557
$ws.Range("L113").Value = 21666.1671
$ws.Range("N113").Value = -26006.1671

$ws.Range("H128").Value = 169854.28
$ws.Range("I128").Value = 169854.28
$ws.Range("K128").Value = 509562.84
$ws.Range("M128").Value = -504582.84

$ws.Range("H137").Value = 3874.7144
$ws.Range("J137").Value = 1661
$ws.Range("L137").Value = 4983
$ws.Range("N137").Value = -15183

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 797.375
$ws.Range("J97").Value = 631.6
$ws.Range("L97").Value = 631.6
$ws.Range("N97").Value = -1623.6

$ws.Range("H102").Value = 8943.036
$ws.Range("I102").Value = 1685.6666
$ws.Range("K102").Value = 1685.6666
$ws.Range("M102").Value = -63.66660000000002

$ws.Range("H126").Value = 7276.5
$ws.Range("I126").Value = 4748.2144
$ws.Range("K126").Value = 14244.6432
$ws.Range("M126").Value = -11774.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2434.9
$ws.Range("J22").Value = 2619.8
$ws.Range("L22").Value = 2619.8
$ws.Range("N22").Value = -3209.8

$ws.Range("H27").Value = 2434.9
$ws.Range("J27").Value = 2619.8
$ws.Range("L27").Value = 2619.8
$ws.Range("N27").Value = -2833.8

$ws.Range("H40").Value = 28497.674
$ws.Range("I40").Value = 51199.773
$ws.Range("J40").Value = 9999.666999999999
$ws.Range("K40").Value = 51199.773
$ws.Range("L40").Value = 9999.666999999999
$ws.Range("M40").Value = -51063.773
$ws.Range("N40").Value = -10271.667

$ws.Range("H82").Value = 1683.625
$ws.Range("J82").Value = 1494.8334
$ws.Range("L82").Value = 1494.8334
$ws.Range("N82").Value = -2216.8334

$ws.Range("H85").Value = 1683.625
$ws.Range("J85").Value = 1494.8334
$ws.Range("L85").Value = 1494.8334
$ws.Range("N85").Value = -3990.8334

$ws.Range("H132").Value = 4262.778
$ws.Range("J132").Value = 9299.833000000001
$ws.Range("L132").Value = 27899.499
$ws.Range("N132").Value = -32959.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 1999
$ws.Range("J19").Value = 1999
$ws.Range("L19").Value = 1999
$ws.Range("N19").Value = -2347

$ws.Range("H30").Value = 10875
$ws.Range("J30").Value = 10875
$ws.Range("L30").Value = 10875
$ws.Range("N30").Value = -11089

$ws.Range("H81").Value = 6470.5713

$ws.Range("H84").Value = 6470.5713

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 4205.8486
$ws.Range("I132").Value = 4028.5356
$ws.Range("K132").Value = 12085.6068
$ws.Range("M132").Value = -9555.606800000001
